# "use fragments for cookieconsent"
#
# Adds two new translation fragments (cookieConsent.message and
# cookieConsent.agree) to the "fragments" sheet, and makes "fragments"
# the active sheet/selection (it was "hidden" before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fragments")

# Write the two new keys first, then their Dutch text values, so the
# shared-string table is populated in the same order as the source
# workbook (keys 467/468, then values 469/470).
$ws.Range("A17").Value = "cookieConsent.message"
$ws.Range("A18").Value = "cookieConsent.agree"
$ws.Range("D17").Value = "Voor deze sites gebruiken we cookies om de gebruikservaring te verbeteren. Indien u verder surft gaan we ervan uit dat u cookies toelaat."
$ws.Range("D18").Value = "Ok"

# Switch the active sheet/selection from "hidden" to "fragments".
$ws.Activate() | Out-Null
$ws.Range("D19").Select() | Out-Null
